$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# Copy the existing "Sprint 2 / PASS" status cell (E10) so the newly
# filled-in status cells pick up the same number-format/border styling
# that the already-completed rows use.
$ws.Range("E10").Copy()
$ws.Range("E11:E19").PasteSpecial(-4122)

# Sprint 2 status column ("PASS") for the user stories that were tested
# this sprint.
$ws.Range("E11").Value = "PASS"
$ws.Range("E12").Value = "PASS"
$ws.Range("E13").Value = "PASS"
$ws.Range("E14").Value = "PASS"
$ws.Range("E15").Value = "PASS"
$ws.Range("E16").Value = "PASS"
$ws.Range("E17").Value = "PASS"
$ws.Range("E18").Value = "PASS"
$ws.Range("E19").Value = "PASS"
$ws.Range("E22").Value = "PASS"

# Sprint 2 tester initials/date comments.
$ws.Range("F11").Value = "JE; 4/2/2018"
$ws.Range("F12").Value = "JE; 4/2/2018"
$ws.Range("F13").Value = "JE; 4/2/2018"
$ws.Range("F14").Value = "JE; 4/2/2018"
$ws.Range("F15").Value = "JE; 4/2/2018"
$ws.Range("F16").Value = "JE; 4/2/2018"
$ws.Range("F17").Value = "JE; 4/2/2018"
$ws.Range("F18").Value = "JE; 4/2/2018"

# Leave the sheet scrolled/selected where the author finished editing.
$ws.Activate()
$ws.Range("F19").Select()
